$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion text (daily rate update) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.47 = 30291.26 pesos`n✅ 30291.26 pesos = 7.42 = 957.74 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 133.9
$wsTasas.Range("O10").Value = 4056
$wsTasas.Range("N12").Value = 4080
